$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.8200883333333334
$ws.Range("H2").Value = 2.460265
$ws.Range("I2").Value = 0.2405117342909232
$ws.Range("J2").Value = 0.2405117342909232
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.352656
$ws.Range("N2").Value = 1.057968
$ws.Range("O2").Value = 0.2868677567428842
$ws.Range("P2").Value = 0.2868677567428842
$ws.Range("Q2").Value = 0.2892090712800001
$ws.Range("R2").Value = 2.60288164152
$ws.Range("S2").Value = 0.06899506168637774
$ws.Range("T2").Value = 0.06899506168637774
$ws.Range("G3").Value = 0.8200883333333334
$ws.Range("H3").Value = 2.460265
$ws.Range("I3").Value = 0.2405117342909232
$ws.Range("J3").Value = 0.2405117342909232
$ws.Range("O3").Value = 0.05437555704326383
$ws.Range("P3").Value = 0.05437555704326383
$ws.Range("Q3").Value = 0.05481935136722223
$ws.Range("R3").Value = 0.493374162305
$ws.Range("S3").Value = 0.01307795952751041
$ws.Range("T3").Value = 0.01307795952751041
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 0.8200883333333334
$ws.Range("H4").Value = 2.460265
$ws.Range("I4").Value = 0.2405117342909232
$ws.Range("J4").Value = 0.2405117342909232
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.050715
$ws.Range("N4").Value = 0.152145
$ws.Range("O4").Value = 0.04125407843114925
$ws.Range("P4").Value = 0.04125407843114925
$ws.Range("Q4").Value = 0.04159077982500001
$ws.Range("R4").Value = 0.374317018425
$ws.Range("S4").Value = 0.009922089950049473
$ws.Range("T4").Value = 0.009922089950049472
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 0.8200883333333334
$ws.Range("H5").Value = 2.460265
$ws.Range("I5").Value = 0.2405117342909232
$ws.Range("J5").Value = 0.2405117342909232
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7591163333333334
$ws.Range("N5").Value = 2.277349
$ws.Range("O5").Value = 0.6175026077827028
$ws.Range("P5").Value = 0.6175026077827028
$ws.Range("Q5").Value = 0.6225424486094445
$ws.Range("R5").Value = 5.602882037485
$ws.Range("S5").Value = 0.1485166231269855
$ws.Range("T5").Value = 0.1485166231269855
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.589676
$ws.Range("H6").Value = 7.769028
$ws.Range("I6").Value = 0.7594882657090768
$ws.Range("J6").Value = 0.7594882657090768
$ws.Range("M6").Value = 0.352656
$ws.Range("N6").Value = 1.057968
$ws.Range("O6").Value = 0.2868677567428842
$ws.Range("P6").Value = 0.2868677567428842
$ws.Range("Q6").Value = 0.913264779456
$ws.Range("R6").Value = 8.219383015104
$ws.Range("S6").Value = 0.2178726950565064
$ws.Range("T6").Value = 0.2178726950565064
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.589676
$ws.Range("H7").Value = 7.769028
$ws.Range("I7").Value = 0.7594882657090768
$ws.Range("J7").Value = 0.7594882657090768
$ws.Range("O7").Value = 0.05437555704326383
$ws.Range("P7").Value = 0.05437555704326383
$ws.Range("Q7").Value = 0.1731086186706667
$ws.Range("R7").Value = 1.557977568036
$ws.Range("S7").Value = 0.04129759751575342
$ws.Range("T7").Value = 0.04129759751575342
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.589676
$ws.Range("H8").Value = 7.769028
$ws.Range("I8").Value = 0.7594882657090768
$ws.Range("J8").Value = 0.7594882657090768
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.050715
$ws.Range("N8").Value = 0.152145
$ws.Range("O8").Value = 0.04125407843114925
$ws.Range("P8").Value = 0.04125407843114925
$ws.Range("Q8").Value = 0.13133541834
$ws.Range("R8").Value = 1.18201876506
$ws.Range("S8").Value = 0.03133198848109978
$ws.Range("T8").Value = 0.03133198848109978
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.589676
$ws.Range("H9").Value = 7.769028
$ws.Range("I9").Value = 0.7594882657090768
$ws.Range("J9").Value = 0.7594882657090768
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.7591163333333334
$ws.Range("N9").Value = 2.277349
$ws.Range("O9").Value = 0.6175026077827028
$ws.Range("P9").Value = 0.6175026077827028
$ws.Range("Q9").Value = 1.965865349641333
$ws.Range("R9").Value = 17.692788146772
$ws.Range("S9").Value = 0.4689859846557172
$ws.Range("T9").Value = 0.4689859846557172